$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 19: Unbreak My Heart | Roof Tile
$ws.Range("H19").Value = 1013.875
$ws.Range("I19").Value = 884.75
$ws.Range("J19").Value = 1143
$ws.Range("K19").Value = 884.75
$ws.Range("L19").Value = 1143
$ws.Range("M19").Value = -709.75
$ws.Range("N19").Value = -1493

# ALC row 28: The Writing Is Not on the Wall | Enchanted Silver Ink
$ws.Range("H28").Value = 467.95834
$ws.Range("I28").Value = 262.77777
$ws.Range("J28").Value = 1083.5
$ws.Range("K28").Value = 262.77777
$ws.Range("L28").Value = 1083.5
$ws.Range("M28").Value = 222.22223
$ws.Range("N28").Value = -2053.5

# ALC row 29: Dripping with Venom | Weak Blinding Potion
$ws.Range("H29").Value = 700
$ws.Range("J29").Value = 1000
$ws.Range("L29").Value = 3000
$ws.Range("N29").Value = -3562

# ALC row 38: Just Give Him a Serum | Hi-Potion of Strength
$ws.Range("H38").Value = 2304250.8
$ws.Range("I38").Value = 2304250.8
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 6912752.399999999
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -6912380.399999999
$ws.Range("N38").ClearContents()

# ALC row 39: Riches' Brew | Hi-Potion of Mind
$ws.Range("H39").Value = 881920.4
$ws.Range("I39").Value = 933795.2
$ws.Range("K39").Value = 2801385.6
$ws.Range("M39").Value = -2801089.6

# ALC row 43: Growing Is Knowing | Growth Formula Gamma
$ws.Range("H43").Value = 5550.5
$ws.Range("I43").Value = 10001
$ws.Range("J43").Value = 1100
$ws.Range("K43").Value = 10001
$ws.Range("L43").Value = 1100
$ws.Range("M43").Value = -9932
$ws.Range("N43").Value = -1238

# ALC row 51: A Bile Business | Shark Oil
$ws.Range("H51").Value = 7863
$ws.Range("J51").Value = 2738.2307
$ws.Range("L51").Value = 2738.2307
$ws.Range("N51").Value = -3706.2307

# ALC row 58: A Matter of Vital Importance | Mega-Potion of Vitality
$ws.Range("H58").Value = 817809.4
$ws.Range("I58").Value = 1032232.9
$ws.Range("K58").Value = 3096698.7
$ws.Range("M58").Value = -3096548.7

# ALC row 64: Forged from the Void | Void Glue
$ws.Range("H64").Value = 58305.11
$ws.Range("I64").Value = 93505.45
$ws.Range("J64").Value = 2990.2856
$ws.Range("K64").Value = 93505.45
$ws.Range("L64").Value = 2990.2856
$ws.Range("M64").Value = -93257.45
$ws.Range("N64").Value = -3486.2856

# ALC row 67: Dodging the Draft (L) | Void Glue
$ws.Range("H67").Value = 58305.11
$ws.Range("I67").Value = 93505.45
$ws.Range("J67").Value = 2990.2856
$ws.Range("K67").Value = 93505.45
$ws.Range("L67").Value = 2990.2856
$ws.Range("M67").Value = -92647.45
$ws.Range("N67").Value = -4706.2856

# ALC row 87: There Was a Late Fee | Noble Gold
$ws.Range("H87").Value = 30522.857
$ws.Range("J87").Value = 30522.857
$ws.Range("L87").Value = 30522.857
$ws.Range("N87").Value = -33018.857

# ALC row 90: A Gate Arcane Is Dragon's Bane (L) | Noble Gold
$ws.Range("H90").Value = 30522.857
$ws.Range("J90").Value = 30522.857
$ws.Range("L90").Value = 91568.571
$ws.Range("N90").Value = -104048.571

# ALC row 98: The Dotted Line | Enchanted Durium Ink
$ws.Range("H98").Value = 456.7143
$ws.Range("I98").Value = 445.36365
$ws.Range("K98").Value = 445.36365
$ws.Range("M98").Value = 1052.63635

# ALC row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Range("H122").Value = 456.7143
$ws.Range("I122").Value = 445.36365
$ws.Range("K122").Value = 1336.09095
$ws.Range("M122").Value = 1113.90905

# ALC row 123: Nearly Bare | Gaja Grimoire
$ws.Range("H123").Value = 38880
$ws.Range("J123").Value = 38880
$ws.Range("L123").Value = 38880
$ws.Range("N123").Value = -48680

# ALC row 129: Practical Command | Commanding Craftsman's Draught
$ws.Range("H129").Value = 910.09576
$ws.Range("I129").Value = 488
$ws.Range("J129").Value = 977.8395
$ws.Range("K129").Value = 1464
$ws.Range("L129").Value = 2933.5185
$ws.Range("M129").Value = 3536
$ws.Range("N129").Value = -12933.5185

# ALC row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 1751.8096
$ws.Range("I137").Value = 1714.25
$ws.Range("K137").Value = 5142.75
$ws.Range("M137").Value = -2592.75

$ws = $wb.Worksheets.Item("ARM")
# ARM row 5: The Alloyed Truth | Bronze Rivets
$ws.Range("H5").Value = 246.75
$ws.Range("I5").Value = 246.75
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 246.75
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -134.75
$ws.Range("N5").ClearContents()

# ARM row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 850.14
$ws.Range("I74").Value = 802.72095
$ws.Range("J74").Value = 1141.4286
$ws.Range("K74").Value = 802.72095
$ws.Range("L74").Value = 1141.4286
$ws.Range("M74").Value = 71.27904999999998
$ws.Range("N74").Value = -2889.4286

# ARM row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 850.14
$ws.Range("I77").Value = 802.72095
$ws.Range("J77").Value = 1141.4286
$ws.Range("K77").Value = 4013.60475
$ws.Range("L77").Value = 5707.143
$ws.Range("M77").Value = 354.39525
$ws.Range("N77").Value = -14443.143

# ARM row 125: The Incomplete Costume | High Durium Armor of Fending
$ws.Range("H125").Value = 48500
$ws.Range("J125").Value = 48500
$ws.Range("L125").Value = 48500
$ws.Range("N125").Value = -58340

$ws = $wb.Worksheets.Item("BSM")
# BSM row 4: Mending Fences | Bronze Rivets
$ws.Range("H4").Value = 246.75
$ws.Range("I4").Value = 246.75
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 246.75
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -131.75
$ws.Range("N4").ClearContents()

# BSM row 80: Unbreaker | Titanium Ingot
$ws.Range("H80").Value = 1146.6666
$ws.Range("I80").Value = 937.5
$ws.Range("J80").Value = 1234.7368
$ws.Range("K80").Value = 937.5
$ws.Range("L80").Value = 1234.7368
$ws.Range("M80").Value = 60.5
$ws.Range("N80").Value = -3230.7368

# BSM row 83: Attack on Titanium (L) | Titanium Ingot
$ws.Range("H83").Value = 1146.6666
$ws.Range("I83").Value = 937.5
$ws.Range("J83").Value = 1234.7368
$ws.Range("K83").Value = 4687.5
$ws.Range("L83").Value = 6173.683999999999
$ws.Range("M83").Value = 304.5
$ws.Range("N83").Value = -16157.684

# BSM row 94: High Steal | High Steel Nugget
$ws.Range("H94").Value = 484.4138
$ws.Range("I94").Value = 454.1875
$ws.Range("J94").Value = 521.61536
$ws.Range("K94").Value = 454.1875
$ws.Range("L94").Value = 521.61536
$ws.Range("M94").Value = -3.1875
$ws.Range("N94").Value = -1423.61536

# BSM row 124: History of the Hrothgar | High Durium Bayonet
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# CRP row 16: Raise the Roof | Ash Lumber
$ws.Range("H16").Value = 619.125
$ws.Range("I16").Value = 615
$ws.Range("J16").Value = 631.5
$ws.Range("K16").Value = 615
$ws.Range("L16").Value = 631.5
$ws.Range("M16").Value = -328
$ws.Range("N16").Value = -1205.5

# CRP row 62: Splinter in the Sewers | Cedar Lumber
$ws.Range("H62").Value = 2394.375
$ws.Range("I62").Value = 1577.5
$ws.Range("J62").Value = 2666.6667
$ws.Range("K62").Value = 1577.5
$ws.Range("L62").Value = 2666.6667
$ws.Range("M62").Value = -953.5
$ws.Range("N62").Value = -3914.6667

# CRP row 65: The Lumber of Their Discontent (L) | Cedar Lumber
$ws.Range("H65").Value = 2394.375
$ws.Range("I65").Value = 1577.5
$ws.Range("J65").Value = 2666.6667
$ws.Range("K65").Value = 7887.5
$ws.Range("L65").Value = 13333.3335
$ws.Range("M65").Value = -4767.5
$ws.Range("N65").Value = -19573.3335

# CRP row 113: Patient Patients | White Ash Lumber
$ws.Range("H113").Value = 619.125
$ws.Range("I113").Value = 615
$ws.Range("J113").Value = 631.5
$ws.Range("K113").Value = 615
$ws.Range("L113").Value = 631.5
$ws.Range("M113").Value = 1555
$ws.Range("N113").Value = -4971.5

# CRP row 124: Earring Awakening | Palm Ear Cuffs of Fending
$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -34910

$ws = $wb.Worksheets.Item("CUL")
# CUL row 14: Keep Your Powder Dry | Kukuru Powder
$ws.Range("H14").Value = 359.8125
$ws.Range("I14").Value = 359.8125
$ws.Range("K14").Value = 1079.4375
$ws.Range("M14").Value = -906.4375

# CUL row 113: Can't Eat Just One | Night Vinegar
$ws.Range("H113").Value = 1847.5
$ws.Range("I113").Value = 3796.6667
$ws.Range("J113").Value = 678
$ws.Range("K113").Value = 11390.0001
$ws.Range("L113").Value = 2034
$ws.Range("M113").Value = -9220.000100000001
$ws.Range("N113").Value = -6374

$ws = $wb.Worksheets.Item("GSM")
# GSM row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 950
$ws.Range("I122").Value = 900
$ws.Range("K122").Value = 2700
$ws.Range("M122").Value = -250

$ws = $wb.Worksheets.Item("LTW")
# LTW row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 4056.923
$ws.Range("I132").Value = 3999.1738
$ws.Range("K132").Value = 11997.5214
$ws.Range("M132").Value = -9467.5214

$ws = $wb.Worksheets.Item("WVR")
# WVR row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 2399.5557
$ws.Range("I122").Value = 1662.25
$ws.Range("J122").Value = 2989.4
$ws.Range("K122").Value = 4986.75
$ws.Range("L122").Value = 8968.200000000001
$ws.Range("M122").Value = -2536.75
$ws.Range("N122").Value = -13868.2
